$wb = $excel.ActiveWorkbook

# Sheet 1: "Overview" -- status columns for zh-cn (E) and de-de (F)
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# Sheet 2: "zh-cn" -- Status column (C)
$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Range("C2").Value = "In Translation"

# Sheet 3: "de-de" -- Status column (C)
$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Range("C2").Value = "In Translation"

# The shorter status text narrows the Status columns; reflect the
# updated (auto-fit-like) column widths on each affected sheet.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
